$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows to append (place, gnd, tgn), matching target row order 14..44
$rows = @(
    @("Czempiń","http://d-nb.info/gnd/10022243-2","http://vocab.getty.edu/tgn/1048994"),
    @("Krotoszyn","http://d-nb.info/gnd/301975-5","http://vocab.getty.edu/tgn/1049186"),
    @("Aalen","http://d-nb.info/gnd/4000015-1","http://vocab.getty.edu/tgn/7005985"),
    @("Amsterdam","http://d-nb.info/gnd/4001783-7","http://vocab.getty.edu/tgn/7006952"),
    @("Baden-Baden","http://d-nb.info/gnd/4004169-4","http://vocab.getty.edu/tgn/7005172"),
    @("Basel","http://d-nb.info/gnd/4004617-5","http://vocab.getty.edu/tgn/7007269"),
    @("Bayreuth","http://d-nb.info/gnd/4005056-7","http://vocab.getty.edu/tgn/7004326"),
    @("Bernburg","http://d-nb.info/gnd/4005774-4","http://vocab.getty.edu/tgn/1036449"),
    @("Bonn","http://d-nb.info/gnd/4007666-0","http://vocab.getty.edu/tgn/7005090"),
    @("Braunau am Inn","http://d-nb.info/gnd/4008037-7","http://vocab.getty.edu/tgn/7003283"),
    @("Bremen","http://d-nb.info/gnd/4008135-7","http://vocab.getty.edu/tgn/7005286"),
    @("Breslau","http://d-nb.info/gnd/4008216-7","http://vocab.getty.edu/tgn/7007601"),
    @("Danzig","http://d-nb.info/gnd/4011039-4","http://vocab.getty.edu/tgn/7007660"),
    @("Darmstadt","http://d-nb.info/gnd/4011077-1","http://vocab.getty.edu/tgn/7004428"),
    @("Davos","http://d-nb.info/gnd/4011220-2","http://vocab.getty.edu/tgn/7007342"),
    @("Dessau","http://d-nb.info/gnd/4011536-7","http://vocab.getty.edu/tgn/7012804"),
    @("Dresden","http://d-nb.info/gnd/4012995-0","http://vocab.getty.edu/tgn/7004455"),
    @("Düsseldorf","http://d-nb.info/gnd/4013255-9","http://vocab.getty.edu/tgn/7004443"),
    @("Eckernförde","http://d-nb.info/gnd/4013498-2","http://vocab.getty.edu/tgn/7012432"),
    @("Eisenach","http://d-nb.info/gnd/4014013-1","http://vocab.getty.edu/tgn/7005969"),
    @("Erlangen","http://d-nb.info/gnd/4015299-6","http://vocab.getty.edu/tgn/7004332"),
    @("Frankfurt am Main","http://d-nb.info/gnd/4018118-2","http://vocab.getty.edu/tgn/7005293"),
    @("Freiburg im Breisgau","http://d-nb.info/gnd/4018272-1","http://vocab.getty.edu/tgn/7004423"),
    @("Genf","http://d-nb.info/gnd/4020137-5","http://vocab.getty.edu/tgn/7007279"),
    @("Gießen","http://d-nb.info/gnd/4020989-1","http://vocab.getty.edu/tgn/7004429"),
    @("Görlitz","http://d-nb.info/gnd/4021441-2","http://vocab.getty.edu/tgn/7012813"),
    @("Göttingen","http://d-nb.info/gnd/4021477-1","http://vocab.getty.edu/tgn/7005246"),
    @("Hechingen","http://d-nb.info/gnd/4023927-5","http://vocab.getty.edu/tgn/1039493"),
    @("Jerusalem","http://d-nb.info/gnd/4028586-8","http://vocab.getty.edu/tgn/7001371"),
    @("Kamen","http://d-nb.info/gnd/4029404-3","http://vocab.getty.edu/tgn/1039810"),
    @("Chemnitz","http://d-nb.info/gnd/4029702-0","http://vocab.getty.edu/tgn/7005951")
)

$startRow = 14
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$lastRow = $startRow + $rows.Count - 1
$nextRow = $lastRow + 1
$ws.Range("A" + $nextRow).Select()
$excel.ActiveWindow.ScrollRow = 22
